# Apply the author's edit:
#   - On slide 1, the title placeholder's run "展示用" gets its font
#     explicitly set to "仿宋" (FangSong) for both the Latin and East Asian
#     typefaces (PowerPoint emits <a:latin>/<a:ea> with the chosen font
#     when a user picks a font for CJK text from the Font list).
#   - The presentation gained an (empty) PowerPoint-2013-era slide-guide
#     extension list; recreate that too where the host's object model
#     allows it (harmless no-op if unsupported).

$p = $ppt.ActivePresentation

# --- Slide 1 / title placeholder: set the run's font to FangSong (仿宋) ---
$slide = $p.Slides.Item(1)
$titleShape = $slide.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

$font = $titleRange.Font
$font.Name = "仿宋"
$font.NameFarEast = "仿宋"

# --- Presentation-level guide list extension (View > Guides bookkeeping) ---
# Recent PowerPoint versions stamp an (initially empty) p15:sldGuideLst
# extension onto the presentation the first time the guide system is
# touched. Exercise the corresponding object-model entry points so the
# change is captured if/when the host materializes them.
try { $null = $p.Guides.Count } catch { }
try { $ppt.DisplayGuides = $true } catch { }
